# A new weekly price record was inserted as row 30, pushing the existing
# rows 30-77 down to rows 31-78 (dimension grows from A1:R77 to A1:R78).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 30; this shifts rows 30-77 down
# to rows 31-78 and carries over the date number format on column D.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new observation.
$ws.Cells.Item(30, 1).Value  = 4
$ws.Cells.Item(30, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value  = "Los Lagos"
$ws.Cells.Item(30, 4).Value  = 44495
$ws.Cells.Item(30, 5).Value  = 10
$ws.Cells.Item(30, 6).Value  = 100112022
$ws.Cells.Item(30, 7).Value  = "Arveja Verde"
$ws.Cells.Item(30, 8).Value  = "Perfection"
$ws.Cells.Item(30, 9).Value  = "Primera"
$ws.Cells.Item(30, 10).Value = 120
$ws.Cells.Item(30, 11).Value = 23000
$ws.Cells.Item(30, 12).Value = 23000
$ws.Cells.Item(30, 13).Value = 23000
$ws.Cells.Item(30, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(30, 16).Value = 920
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
